$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabula-99res_2")

# Move the workbook window horizontally (xWindow: -108 -> 22932), matching the
# author's on-screen window position when the file was re-saved.
$win = $excel.ActiveWindow
$win.Left = 22932
$win.Top = -108

# Update "Total new nominations" value (B32): 37294 -> 36294
$ws.Range("B32").Value = 36294

# Update "Total confirmed" value (B34): 39983 -> 39893
$ws.Range("B34").Value = 39893
